$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace the HFR / 791621.4 record with the BERTILLA / 50442.09 record.
# Force column A to Text format first so the leading zeros in the account
# number ("005064906") are preserved instead of being coerced to a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "005064906"
$ws.Range("B2").Value = "BERTILLA"
$ws.Range("C2").Value = 50442.09

# The old BERTILLA / 441.97 row (further down the sheet, at row 41) is removed
# entirely now that its data has moved up to row 2; delete the whole row so
# everything below shifts up by one.
$ws.Rows("41").Delete()
